$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1400
$ws.Range("B3").Value = 1176
$ws.Range("B4").Value = 1512
$ws.Range("B5").Value = 1064
$ws.Range("B6").Value = 1064
$ws.Range("B7").Value = 1136
$ws.Range("B8").Value = 1504
$ws.Range("B9").Value = 1008
$ws.Range("B10").Value = 361

$ws.Range("B2:B10").Select()
